$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Burndown chart update: Day 20 (column V) time logged for the two user
# stories (rows 2 and 3).

# US#1 Class Amis (Nicolas Max) - Day 20 already uses the "hours logged"
# highlighted format; just update the value.
$ws.Range("V2").Value = 2

# US#2 Class User (Eliott Vincenzo) - Day 20 gets the same highlighted
# format as V2 (copy formatting) plus the new value.
$ws.Range("V2").Copy()
$ws.Range("V3").PasteSpecial(-4122)
$ws.Range("V3").Value = 2

# Cosmetic view state: zoom level and active cell selection as left by
# the author.
$excel.ActiveWindow.Zoom = 85
$ws.Range("W6").Select() | Out-Null
